# Weekly CompStat update: new crime data collected (week of 2/13/2023 - 2/19/2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume/number + reporting week dates ---
$ws.Range("A8").Value = "Volume 30   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/13/2023  Through  2/19/2023"

# --- Anchor cells used as Copy() sources so text-placeholder cells ("0" / "***.*") keep their
#     original shared-string + style (rather than becoming new number-formatted strings). ---
$zeroSrc = $ws.Range("C14")   # literal text "0", style 14
$naSrc   = $ws.Range("E14")   # literal text "***.*", style 14

# Row 15 - Rape
$zeroSrc.Copy($ws.Range("C15"))
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 200
$ws.Range("N15").Value = -80

# Row 16 - Robbery
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -16.666666666666
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -14.285714285714
$ws.Range("I16").Value = 21
$ws.Range("J16").Value = 21
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = -25
$ws.Range("M16").Value = -51.162790697674
$ws.Range("N16").Value = -92.881355932203

# Row 17 - Fel. Assault
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = -20.689655172413
$ws.Range("I17").Value = 48
$ws.Range("J17").Value = 56
$ws.Range("K17").Value = -14.285714285714
$ws.Range("L17").Value = -2.04081632653
$ws.Range("M17").Value = 4.347826086956
$ws.Range("N17").Value = -64.705882352941

# Row 18 - Burglary
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("I18").Value = 24
$ws.Range("J18").Value = 23
$ws.Range("K18").Value = 4.347826086956
$ws.Range("L18").Value = -17.241379310344
$ws.Range("M18").Value = -54.716981132075
$ws.Range("N18").Value = -86.206896551724

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -12.5
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 25
$ws.Range("H19").Value = 28
$ws.Range("I19").Value = 55
$ws.Range("J19").Value = 50
$ws.Range("K19").Value = 10
$ws.Range("L19").Value = 41.025641025641
$ws.Range("M19").Value = 25
$ws.Range("N19").Value = -46.078431372549

# Row 20 - G.L.A.
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 80
$ws.Range("I20").Value = 14
$ws.Range("J20").Value = 12
$ws.Range("K20").Value = 16.666666666666
$ws.Range("L20").Value = 40
$ws.Range("M20").Value = -17.647058823529
$ws.Range("N20").Value = -84.782608695652

# Row 21 - TOTAL
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -15.384615384615
$ws.Range("F21").Value = 95
$ws.Range("H21").Value = 7.954545454545
$ws.Range("I21").Value = 165
$ws.Range("J21").Value = 167
$ws.Range("K21").Value = -1.197604790419
$ws.Range("L21").Value = 3.77358490566
$ws.Range("M21").Value = -20.289855072463
$ws.Range("N21").Value = -79.804161566707

# Row 22 - Transit
$ws.Range("G22").Value = 1

# Row 23 - Housing
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 15
$ws.Range("H23").Value = -33.333333333333
$ws.Range("I23").Value = 28
$ws.Range("J23").Value = 31
$ws.Range("K23").Value = -9.677419354838
$ws.Range("L23").Value = -9.677419354838
$ws.Range("M23").Value = 75

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -29.032258064516
$ws.Range("F24").Value = 109
$ws.Range("G24").Value = 124
$ws.Range("H24").Value = -12.096774193548
$ws.Range("I24").Value = 200
$ws.Range("J24").Value = 197
$ws.Range("K24").Value = 1.522842639593
$ws.Range("L24").Value = 43.88489208633
$ws.Range("M24").Value = 63.934426229508

# Row 25 - Misd. Assault
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 42
$ws.Range("G25").Value = 51
$ws.Range("H25").Value = -17.647058823529
$ws.Range("I25").Value = 91
$ws.Range("J25").Value = 92
$ws.Range("K25").Value = -1.086956521739
$ws.Range("L25").Value = 65.454545454545
$ws.Range("M25").Value = -5.208333333333

# Row 26 - UCR Rape*
$zeroSrc.Copy($ws.Range("C26"))
$zeroSrc.Copy($ws.Range("D26"))
$naSrc.Copy($ws.Range("E26"))
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 50

# Row 27 - Other Sex Crimes
$ws.Range("C27").Value = 1
$ws.Range("I27").Value = 5
$ws.Range("K27").Value = 66.666666666666
$ws.Range("L27").Value = 150

# Row 28 - Shooting Vic.
$zeroSrc.Copy($ws.Range("G28"))
$naSrc.Copy($ws.Range("H28"))

# Row 29 - Shooting Inc.
$zeroSrc.Copy($ws.Range("G29"))
$naSrc.Copy($ws.Range("H29"))
